$d = $word.ActiveDocument
$lastPara = $d.Paragraphs.Last
$origStart = $lastPara.Range.Start
$origEnd = $lastPara.Range.End

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:widowControl/><w:bidi w:val="0"/><w:ind w:hanging="0" w:start="0" w:end="0"/><w:jc w:val="start"/><w:rPr><w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/><w:sz w:val="24"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>        <w:t xml:space="preserve">স্ট্রিংকে আরও সাধারণ অর্থে বলা হয় ব্রেইন। বহুমাত্রিক মেমব্রেইন থেকে নামটা আসা। এরা এত ক্ষুদ্র যে কোনো যন্ত্র দিয়েই এদেরকে দেখার কোন আশা নেই। অন্তত আমাদের সভ্যতার আরও অনেক অনেক বেশি উন্নত না হওয়া পর্যন্ত অপেক্ষা করতে হবেই। কনাপদার্থবিদরা অতিপারমাণবিক জগত দেখেন কণাত্বরকের </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="SolaimanLipi" w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>(</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>যে যন্ত্র কণাকে বিশাল বেগে চালিত করে</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="SolaimanLipi" w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">) </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>সাহায্যে। চুম্বকক্ষেত্র বা অন্যকোনো ক্ষেত্র কাজে লাগিয়ে ক্ষুদ্র ক্ষুদ্র কণাকে অনেক বেশি বেগে ধাবিত করা হয়।এ কণারা একে অপরের সাথে সংঘর্ষ করে বিভিন্ন অংশ আলাদা হয়ে যায়। কনাত্বুরক হলো অতিপারমাণবিক জগতের মাইক্রোস্কোপ। এসব কণায় যত বেশি শক্তি দেওয়া হবে</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="SolaimanLipi" w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">মাইক্রোস্কোপ তত শক্তিশালী হবে। আর ততই ক্ষুদ্র বস্তু দেখা যাবে। </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:widowControl/>
        <w:bidi w:val="0"/>
        <w:ind w:hanging="0" w:start="0" w:end="0"/>
        <w:jc w:val="start"/>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:widowControl/>
        <w:bidi w:val="0"/>
        <w:ind w:hanging="0" w:start="0" w:end="0"/>
        <w:jc w:val="start"/>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>এমন এক কণাত্বরকের নাম সুপাকন্ডাক্টিং সুপার কোলাইডার। ১৯৯০</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="SolaimanLipi" w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>-</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>এর দশকের শুরু পর্যন্ত বহু বিলিয়ন ডলারের এ যন্ত্র বানানোর পরিকল্পনা ছিল। নির্মিত হলে এটি হত সর্বকালের সবচেয়ে শক্তিশালী কণাত্বরক। ৫৪ মাইল লম্বা লুপের মধ্যে ১০</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="SolaimanLipi" w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>০০০ চুম্বক থাকার কথা ছিল এতে। এত শক্তিশালী যন্ত্রও গুটিয়ে থাকা মাত্রাগুলো বা স্ট্রিং দেখার জন্য যথেষ্ট নয়। হ্যাঁ</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="SolaimanLipi" w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>কণাত্বুরক দিয়ে স্ট্রিং দেখার ব্যবস্থা করা যাবে। সেজন্য কণাত্বরককে হতে হবে ৬০ কোটি</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="SolaimanLipi" w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>-</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>কোটি মাইলের একটি লুপ। একটি কণা আলোর বেগে চললেও এত বড় দূরত্ব পাড়ি দিতে সময় লাগবে ১</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="SolaimanLipi" w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">০০০ বছর। </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:widowControl/>
        <w:bidi w:val="0"/>
        <w:ind w:hanging="0" w:start="0" w:end="0"/>
        <w:jc w:val="start"/>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:widowControl/>
        <w:bidi w:val="0"/>
        <w:ind w:hanging="0" w:start="0" w:end="0"/>
        <w:jc w:val="start"/>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>বর্তমানে কল্পনাযোগ্য কোনো যন্ত্র দিয়েই বিজ্ঞানীরা সরাসরি স্ট্রিং দেখার আশা করতে পারেন না। এমন কোনো পরীক্ষার সন্ধান কেউ দিতে পারবে না</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="SolaimanLipi" w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>যার মাধ্যমে যাচাই করে দেখা যাবে ব্ল্যাকহোল ও কণা আসলেই স্ট্রিং কিনা। স্ট্রিং তত্ত্বের বিরুদ্ধে এটাই সবচেয়ে বড় আপত্তি। বিজ্ঞান হলো পর্যবেক্ষণ ও পরীক্ষা নির্ভর। তাই কেউ কেউ বলেন</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="SolaimanLipi" w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>স্ট্রিং তত্ত্ব আসলে বিজ্ঞান নয়</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="SolaimanLipi" w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">বরং দর্শন। </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="SolaimanLipi" w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>(</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>সাম্প্রতিক এক গুচ্ছ তত্ত্ব বলছে</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="SolaimanLipi" w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>কিছু সংকীর্ণ মাত্রা লম্বায় ১০</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="SolaimanLipi" w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:vertAlign w:val="superscript"/>
        </w:rPr>
        <w:t>-</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:vertAlign w:val="superscript"/>
        </w:rPr>
        <w:t>১৯</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> সেন্টিমিটার বা আরও বড় হতে পারে। সেক্ষেত্রে এগুলো পরীক্ষারযোগ্য হবে। তবে এখনও পর্যন্ত এ তত্ত্বগুলোকে ভুলই মনে হচ্ছে। ভাবনাগুলো চমৎকার। কিন্তু সঠিক হওয়ার সম্ভাবনা সামান্য। </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:widowControl/>
        <w:bidi w:val="0"/>
        <w:ind w:hanging="0" w:start="0" w:end="0"/>
        <w:jc w:val="start"/>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:widowControl/>
        <w:bidi w:val="0"/>
        <w:ind w:hanging="0" w:start="0" w:end="0"/>
        <w:jc w:val="start"/>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">নিউটন মহাকর্ষ ও গতিসূত্র আবিষ্কার করেছিলেন। এর মাধ্যমে মহাবিশ্বের মধ্য দিয়ে গ্রহ ও অন্যান্য বস্তুর চলাচলের ব্যাখ্যা পেয়েছিলেন পদার্থবিদরা। নতুন ধূমকেতু আবিষ্কৃত হলেই নিউটনের হিসাবের পক্ষে সমর্থন জোরালো হত। তবে কিছু সমস্যা ছিল। একটি সমস্যা ছিল বুধ গ্রহের কক্ষপথ। গ্রহটির কক্ষপথ যেভাবে দুলে ওঠে তা মেলে না নিউটনের পূর্বাভাসের সাথে। তবে সার্বিকভাবে নিউটনের তত্ত্ব একের পর একে পরীক্ষায় পাশ করে যাচ্ছিল। </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:widowControl/>
        <w:bidi w:val="0"/>
        <w:ind w:hanging="0" w:start="0" w:end="0"/>
        <w:jc w:val="start"/>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:widowControl/>
        <w:bidi w:val="0"/>
        <w:ind w:hanging="0" w:start="0" w:end="0"/>
        <w:jc w:val="start"/>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">আইনস্টাইনের তত্ত্ব নিউটনের ভুলগুলো সংশোধন করে। সমাধান হত বুধ গ্রহের গতির। এ তত্ত্বগুলো মহাকর্ষ সম্পর্কে পরীক্ষাযোগ্য পূর্বানুমানও করে। সূর্যগ্রহণের সময় এডিংটন নক্ষত্রের আলোর বাঁক পর্যবেক্ষণ করেন। সত্য হয় একটি অনুমান। </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:widowControl/>
        <w:bidi w:val="0"/>
        <w:ind w:hanging="0" w:start="0" w:end="0"/>
        <w:jc w:val="start"/>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:widowControl/>
        <w:bidi w:val="0"/>
        <w:ind w:hanging="0" w:start="0" w:end="0"/>
        <w:jc w:val="start"/>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>অন্যদিকে স্ট্রিং তত্ত্ব অনেকগুলো প্রচলিত তত্ত্বকে সুন্দরভাবে জোড়া দেয়। ব্ল্যাকহোল ও কণার আচরণ সম্পর্কে বেশ কিছু পূর্বাভাসও দেয়। তবে এগুলোর কোনোটিই পরীক্ষাযোগ্য নয়। নয় পর্যবেক্ষণযোগ্য। গাণিতিকভাবে হয়ত তত্ত্বটা সুসঙ্গত। এমনকি সুন্দরও। তবে এটা এখনও বিজ্ঞানের অংশ নয়</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="32"/>
          <w:vertAlign w:val="superscript"/>
        </w:rPr>
        <w:t>১</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">। </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:widowControl/>
        <w:bidi w:val="0"/>
        <w:ind w:hanging="0" w:start="0" w:end="0"/>
        <w:jc w:val="start"/>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:widowControl/>
        <w:bidi w:val="0"/>
        <w:ind w:hanging="0" w:start="0" w:end="0"/>
        <w:jc w:val="start"/>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">তথ্যনির্দেশ </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:widowControl/>
        <w:bidi w:val="0"/>
        <w:ind w:hanging="0" w:start="0" w:end="0"/>
        <w:jc w:val="start"/>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:widowControl/>
        <w:bidi w:val="0"/>
        <w:ind w:hanging="0" w:start="0" w:end="0"/>
        <w:jc w:val="start"/>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="SolaimanLipi" w:hAnsi="SolaimanLipi" w:cs="SolaimanLipi"/>
          <w:sz w:val="24"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">১।                </w:t>

      </w:r>
    </w:p>
'@

$insertionPoint = $d.Range($origEnd, $origEnd)
$insertionPoint.InsertXML($xml)

$oldRange = $d.Range($origStart, $origEnd)
$oldRange.Delete()

Write-Output "done"
